$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 177, shifting the existing row 177 (and all
# rows below it) down by one. This matches the diff, where a new data
# record (date 2021-11-08 / 44508) was inserted ahead of the previously
# first "Provincia de Cautín" record, pushing every subsequent row down by
# one position (dimension grows from A1:R216 to A1:R217).
$ws.Rows.Item(177).Insert()

# Populate the newly inserted row 177 with its data. The categorical
# columns (A, B, C, E, F, G, H, I, N, O, Q, R) mirror the record that used
# to sit in row 177, while D/J/K/L/M/P carry the new values from the diff.
$ws.Range("A177").Value = 10
$ws.Range("B177").Value = "Vega Modelo de Temuco"
$ws.Range("C177").Value = "La Araucanía"
$ws.Range("D177").Value = 44508
$ws.Range("E177").Value = 9
$ws.Range("F177").Value = 100112009
$ws.Range("G177").Value = "Acelga"
$ws.Range("H177").Value = "Sin especificar"
$ws.Range("I177").Value = "Primera"
$ws.Range("J177").Value = 50
$ws.Range("K177").Value = 8000
$ws.Range("L177").Value = 9000
$ws.Range("M177").Value = 8600
$ws.Range("N177").Value = "$/docena de atados (12 kilos)"
$ws.Range("O177").Value = "Provincia de Cautín"
$ws.Range("P177").Value = 717
$ws.Range("Q177").Value = 12
$ws.Range("R177").Value = "Hortaliza"
